$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels for BP1/BQ1 (average_doctor <-> average_doctor_old)
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Updated stats values (recomputed averages incl. new doctor case)
$ws.Range("E4").Value = 0.485
$ws.Range("F4").Value = 0.055
$ws.Range("G4").Value = 0.235
$ws.Range("N4").Value = 0.482
$ws.Range("O4").Value = 0.057
$ws.Range("P4").Value = 0.24
$ws.Range("Q4").Value = 0.035
$ws.Range("R4").Value = 0.028
$ws.Range("S4").Value = 0.167
$ws.Range("W4").Value = 0.353
$ws.Range("X4").Value = 0.104
$ws.Range("Y4").Value = 0.323
$ws.Range("AI4").Value = 0.425
$ws.Range("AJ4").Value = 0.096
$ws.Range("AK4").Value = 0.309
$ws.Range("AU4").Value = 0.234
$ws.Range("AV4").Value = 0.028
$ws.Range("AW4").Value = 0.167
$ws.Range("BA4").Value = 2.048
$ws.Range("BB4").Value = 0.137
$ws.Range("BC4").Value = 0.369
$ws.Range("BG4").Value = 0.728
$ws.Range("BH4").Value = 0.14
$ws.Range("BI4").Value = 0.374
$ws.Range("BM4").Value = 0.75
$ws.Range("BN4").Value = 0.064
$ws.Range("BO4").Value = 0.253
$ws.Range("BP4").Value = 0.6830000000000001
$ws.Range("BQ4").Value = 0.756
$ws.Range("E5").Value = 0.612
$ws.Range("F5").Value = 0.063
$ws.Range("G5").Value = 0.252
$ws.Range("N5").Value = 0.735
$ws.Range("O5").Value = 0.081
$ws.Range("P5").Value = 0.285
$ws.Range("Q5").Value = 0.014
$ws.Range("R5").Value = 0.004
$ws.Range("S5").Value = 0.061
$ws.Range("W5").Value = 0.342
$ws.Range("X5").Value = 0.106
$ws.Range("Y5").Value = 0.326
$ws.Range("AI5").Value = 0.426
$ws.Range("AJ5").Value = 0.091
$ws.Range("AK5").Value = 0.302
$ws.Range("AU5").Value = 0.436
$ws.Range("AV5").Value = 0.083
$ws.Range("AW5").Value = 0.289
$ws.Range("BA5").Value = 1.33
$ws.Range("BG5").Value = 0.384
$ws.Range("BH5").Value = 0.051
$ws.Range("BI5").Value = 0.226
$ws.Range("BM5").Value = 0.539
$ws.Range("BN5").Value = 0.049
$ws.Range("BO5").Value = 0.221
$ws.Range("BP5").Value = 0.443
$ws.Range("BQ5").Value = 0.46
$ws.Range("E6").Value = 0.541
$ws.Range("N6").Value = 0.582
$ws.Range("Q6").Value = 0.02
$ws.Range("W6").Value = 0.347
$ws.Range("AI6").Value = 0.425
$ws.Range("AU6").Value = 0.305
$ws.Range("BA6").Value = 1.605
$ws.Range("BG6").Value = 0.503
$ws.Range("BM6").Value = 0.627
$ws.Range("BP6").Value = 0.535
$ws.Range("BQ6").Value = 0.569
$ws.Range("E7").Value = 0.582
$ws.Range("N7").Value = 0.665
$ws.Range("Q7").Value = 0.016
$ws.Range("W7").Value = 0.344
$ws.Range("AI7").Value = 0.426
$ws.Range("AU7").Value = 0.372
$ws.Range("BA7").Value = 1.427
$ws.Range("BG7").Value = 0.424
$ws.Range("BM7").Value = 0.571
$ws.Range("BP7").Value = 0.476
$ws.Range("BQ7").Value = 0.498
$ws.Range("E8").Value = 0.702
$ws.Range("F8").Value = 0.074
$ws.Range("G8").Value = 0.272
$ws.Range("N8").Value = 0.824
$ws.Range("O8").Value = 0.052
$ws.Range("P8").Value = 0.229
$ws.Range("Q8").Value = 0.021
$ws.Range("S8").Value = 0.11
$ws.Range("W8").Value = 0.386
$ws.Range("X8").Value = 0.121
$ws.Range("Y8").Value = 0.348
$ws.Range("AI8").Value = 0.491
$ws.Range("AJ8").Value = 0.133
$ws.Range("AK8").Value = 0.365
$ws.Range("AU8").Value = 0.375
$ws.Range("AV8").Value = 0.09
$ws.Range("AW8").Value = 0.3
$ws.Range("BA8").Value = 1.773
$ws.Range("BB8").Value = 0.106
$ws.Range("BC8").Value = 0.326
$ws.Range("BG8").Value = 0.5610000000000001
$ws.Range("BH8").Value = 0.11
$ws.Range("BI8").Value = 0.332
$ws.Range("BM8").Value = 0.677
$ws.Range("BN8").Value = 0.062
$ws.Range("BO8").Value = 0.248
$ws.Range("BP8").Value = 0.591
$ws.Range("BQ8").Value = 0.622
$ws.Range("E9").Value = 0.658
$ws.Range("F9").Value = 0.225
$ws.Range("G9").Value = 0.474
$ws.Range("N9").Value = 0.763
$ws.Range("O9").Value = 0.181
$ws.Range("P9").Value = 0.425
$ws.Range("W9").Value = 0.263
$ws.Range("X9").Value = 0.194
$ws.Range("Y9").Value = 0.44
$ws.Range("AI9").Value = 0.421
$ws.Range("AJ9").Value = 0.244
$ws.Range("AK9").Value = 0.494
$ws.Range("BA9").Value = 1.711
$ws.Range("BB9").Value = 0.249
$ws.Range("BC9").Value = 0.499
$ws.Range("BG9").Value = 0.605
$ws.Range("BH9").Value = 0.239
$ws.Range("BI9").Value = 0.489
$ws.Range("BM9").Value = 0.632
$ws.Range("BN9").Value = 0.233
$ws.Range("BO9").Value = 0.482
$ws.Range("BP9").Value = 0.57
$ws.Range("BQ9").Value = 0.607
$ws.Range("E10").Value = 0.8159999999999999
$ws.Range("F10").Value = 0.15
$ws.Range("G10").Value = 0.388
$ws.Range("N10").Value = 0.947
$ws.Range("O10").Value = 0.05
$ws.Range("P10").Value = 0.223
$ws.Range("W10").Value = 0.474
$ws.Range("X10").Value = 0.249
$ws.Range("Y10").Value = 0.499
$ws.Range("AI10").Value = 0.526
$ws.Range("AJ10").Value = 0.249
$ws.Range("AK10").Value = 0.499
$ws.Range("AU10").Value = 0.368
$ws.Range("AV10").Value = 0.233
$ws.Range("AW10").Value = 0.482
$ws.Range("BA10").Value = 2.185
$ws.Range("BB10").Value = 0.206
$ws.Range("BC10").Value = 0.454
$ws.Range("BG10").Value = 0.658
$ws.Range("BH10").Value = 0.225
$ws.Range("BI10").Value = 0.474
$ws.Range("BM10").Value = 0.8159999999999999
$ws.Range("BN10").Value = 0.15
$ws.Range("BO10").Value = 0.388
$ws.Range("BP10").Value = 0.728
$ws.Range("BQ10").Value = 0.76
$ws.Range("E11").Value = 0.842
$ws.Range("F11").Value = 0.133
$ws.Range("G11").Value = 0.365
$ws.Range("N11").Value = 0.947
$ws.Range("O11").Value = 0.05
$ws.Range("P11").Value = 0.223
$ws.Range("W11").Value = 0.474
$ws.Range("X11").Value = 0.249
$ws.Range("Y11").Value = 0.499
$ws.Range("AI11").Value = 0.605
$ws.Range("AJ11").Value = 0.239
$ws.Range("AK11").Value = 0.489
$ws.Range("AU11").Value = 0.5
$ws.Range("AV11").Value = 0.25
$ws.Range("AW11").Value = 0.5
$ws.Range("BA11").Value = 2.185
$ws.Range("BB11").Value = 0.206
$ws.Range("BC11").Value = 0.454
$ws.Range("BG11").Value = 0.658
$ws.Range("BH11").Value = 0.225
$ws.Range("BI11").Value = 0.474
$ws.Range("BM11").Value = 0.8159999999999999
$ws.Range("BN11").Value = 0.15
$ws.Range("BO11").Value = 0.388
$ws.Range("BP11").Value = 0.728
$ws.Range("BQ11").Value = 0.76
$ws.Range("E12").Value = 1.344
$ws.Range("F12").Value = 0.538
$ws.Range("G12").Value = 0.734
$ws.Range("N12").Value = 1.25
$ws.Range("O12").Value = 0.299
$ws.Range("P12").Value = 0.546
$ws.Range("W12").Value = 1.556
$ws.Range("X12").Value = 0.469
$ws.Range("Y12").Value = 0.6850000000000001
$ws.Range("AI12").Value = 1.652
$ws.Range("AJ12").Value = 1.531
$ws.Range("AK12").Value = 1.237
$ws.Range("AU12").Value = 2.857
$ws.Range("AV12").Value = 3.646
$ws.Range("AW12").Value = 1.91
$ws.Range("BA12").Value = 3.782
$ws.Range("BB12").Value = 0.469
$ws.Range("BC12").Value = 0.6850000000000001
$ws.Range("BG12").Value = 1.08
$ws.Range("BH12").Value = 0.074
$ws.Range("BI12").Value = 0.271
$ws.Range("BM12").Value = 1.258
$ws.Range("BN12").Value = 0.256
$ws.Range("BO12").Value = 0.506
$ws.Range("BP12").Value = 1.261
$ws.Range("BQ12").Value = 1.229
$ws.Range("E13").Value = 1.421
$ws.Range("F13").Value = 0.319
$ws.Range("G13").Value = 0.5649999999999999
$ws.Range("N13").Value = 1.714
$ws.Range("O13").Value = 0.438
$ws.Range("P13").Value = 0.662
$ws.Range("W13").Value = 1.025
$ws.Range("X13").Value = 0.185
$ws.Range("Y13").Value = 0.43
$ws.Range("AI13").Value = 1.164
$ws.Range("AJ13").Value = 0.32
$ws.Range("AK13").Value = 0.5649999999999999
$ws.Range("AU13").Value = 2.056
$ws.Range("AV13").Value = 0.349
$ws.Range("AW13").Value = 0.591
$ws.Range("BA13").Value = 2.196
$ws.Range("BB13").Value = 0.286
$ws.Range("BC13").Value = 0.535
$ws.Range("BG13").Value = 0.541
$ws.Range("BH13").Value = 0.054
$ws.Range("BI13").Value = 0.232
$ws.Range("BM13").Value = 0.806
$ws.Range("BN13").Value = 0.168
$ws.Range("BO13").Value = 0.41
$ws.Range("BP13").Value = 0.732
$ws.Range("BQ13").Value = 0.671
